# Update the "Test Data" sheet: clear the contents of row 24 (A24:H24),
# leaving the existing number formatting (styles) on A24/D24/E24/F24 intact,
# and move the selection to the now-empty row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")

$ws.Activate()
$ws.Range("A24:H24").ClearContents()
$ws.Range("A24:H24").Select()
